# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: four long space-separated simulation-distribution strings
# get extra numbers appended (one new week's worth of sim results).
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$old = $ydsWs.Range("B2").Value2
$ydsWs.Range("B2").Value = "$old 5 1 4 5 4 3 -4 -1 7 2 18 4 8 5 0 11 13 21 5 3 2 1 1"

$old = $ydsWs.Range("B3").Value2
$ydsWs.Range("B3").Value = "$old 11 2 -1 5 3 5 4 9 4 20 10 32 4 8 -3 20 36 12 3 11 2 6 5 6 4 8 14"

$old = $ydsWs.Range("C2").Value2
$ydsWs.Range("C2").Value = "$old 2 1 3 3 2 5 1 6 15 1 2 13 -2 5 6 3 4 11 8 11 1 2 7 1 0 -2 1"

$old = $ydsWs.Range("C3").Value2
$ydsWs.Range("C3").Value = "$old -1 8 3 9 4 34 5 5 9 7 6 27 9 1 12 5 7 13 17 3 3 4"

# ---------------------------------------------------------------------
# OFF sheet: running totals updated with week 15 numbers.
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 11
$offWs.Range("C2").Value = 444
$offWs.Range("D2").Value = 22
$offWs.Range("E2").Value = 25
$offWs.Range("F2").Value = 131
$offWs.Range("G2").Value = 134
$offWs.Range("I2").Value = 26
$offWs.Range("J2").Value = 108
$offWs.Range("N2").Value = 59
$offWs.Range("O2").Value = 38
$offWs.Range("P2").Value = 24

$offWs.Range("C3").Value = 302
$offWs.Range("E3").Value = 58
$offWs.Range("F3").Value = 201
$offWs.Range("G3").Value = 41
$offWs.Range("H3").Value = 55
$offWs.Range("I3").Value = 108
$offWs.Range("J3").Value = 57
$offWs.Range("L3").Value = 500
$offWs.Range("M3").Value = 324
$offWs.Range("Q3").Value = 1106

# ---------------------------------------------------------------------
# DEF sheet: same shape of running totals.
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 353
$defWs.Range("D2").Value = 21
$defWs.Range("E2").Value = 14
$defWs.Range("F2").Value = 90
$defWs.Range("G2").Value = 85
$defWs.Range("H2").Value = 6
$defWs.Range("I2").Value = 12
$defWs.Range("J2").Value = 39

$defWs.Range("B3").Value = 19
$defWs.Range("C3").Value = 336
$defWs.Range("D3").Value = 12
$defWs.Range("F3").Value = 214
$defWs.Range("G3").Value = 62
$defWs.Range("H3").Value = 55
$defWs.Range("I3").Value = 114
$defWs.Range("J3").Value = 98
$defWs.Range("L3").Value = 617
$defWs.Range("M3").Value = 380
$defWs.Range("Q3").Value = 1031

# ---------------------------------------------------------------------
# ST sheet: numeric totals plus the six distribution strings.
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 170
$stWs.Range("D2").Value = 116
$stWs.Range("F2").Value = 378
$stWs.Range("G2").Value = 374
$stWs.Range("H2").Value = 9
$stWs.Range("L2").Value = 108
$stWs.Range("M2").Value = 97
$stWs.Range("N2").Value = 65
$stWs.Range("O2").Value = 47
$stWs.Range("B3").Value = 109

$old = $stWs.Range("B4").Value2
$stWs.Range("B4").Value = "$old 63 63 47 40"

$old = $stWs.Range("B5").Value2
$stWs.Range("B5").Value = "$old 21 0 13 6"

$old = $stWs.Range("B6").Value2
$stWs.Range("B6").Value = "$old 19 12 23"

$old = $stWs.Range("D3").Value2
$stWs.Range("D3").Value = "$old 60 44 44"

$old = $stWs.Range("D4").Value2
$stWs.Range("D4").Value = "$old 8 1 6"

$old = $stWs.Range("D5").Value2
$stWs.Range("D5").Value = "$old 0 0 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet.
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C3").Value = 11
$turnsWs.Range("D3").Value = 21

# ---------------------------------------------------------------------
# PEN sheet.
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B3").Value = 44
$penWs.Range("B4").Value = 2
$penWs.Range("D4").Value = 25
$penWs.Range("B5").Value = 2
